$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The columns D (mes-nombre), I (sector-descripcion) and J (sexo) are
# re-curated from "dimension" fields into "measure" fields.

# Row 2: iaest-dimension:<x> -> iaest-measure:<x>
$ws.Range("D2").Value = "iaest-measure:mes-nombre"
$ws.Range("I2").Value = "iaest-measure:sector-descripcion"
$ws.Range("J2").Value = "iaest-measure:sexo"

# Row 3: dim -> medida
$ws.Range("D3").Value = "medida"
$ws.Range("I3").Value = "medida"
$ws.Range("J3").Value = "medida"

# Row 4: skos:Concept -> xsd:int
$ws.Range("D4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("J4").Value = "xsd:int"

# Row 5: measures no longer need a mapping file, so the corresponding
# cells are removed entirely (not just blanked).
$ws.Range("D5").Clear()
$ws.Range("I5").Clear()
$ws.Range("J5").Clear()

"done"
